$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 5 = CU-01 : Esfuerzo (F5) 2 -> 1
$ws.Range("F5").Value = 1

# Row 6 = CU-02 : Estado (E6) vacio -> planificado ; Esfuerzo (F6) 0 -> 1
$ws.Range("E6").Value = "planificado"
$ws.Range("F6").Value = 1

# Row 7 = CU-03 : Estado (E7) vacio -> planificado ; Esfuerzo (F7) 0 -> 1
$ws.Range("E7").Value = "planificado"
$ws.Range("F7").Value = 1

# Row 8 = CU-04 : Esfuerzo (F8) 2 -> 1
$ws.Range("F8").Value = 1

# Row 9 = CU-05 : Estado (E9) vacio -> planificado ; Esfuerzo (F9) 0 -> 1
$ws.Range("E9").Value = "planificado"
$ws.Range("F9").Value = 1

# Row 11 = CU-07 : Estado (E11) vacio -> planificado ; Esfuerzo (F11) 0 -> 1
$ws.Range("E11").Value = "planificado"
$ws.Range("F11").Value = 1

# Update the active selection on the sheet to E12 as in the diff
$ws.Range("E12").Select()
